$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-obsolete rows 9-19 first (data set shrinks to A1:D8)
$ws.Range("A9:D19").ClearContents()

# Force text format on columns A, C, D so numeric-looking strings stay as text
# (matches the workbook's existing convention of storing these lists as text)
$ws.Range("A2:A8").NumberFormat = "@"
$ws.Range("C2:D8").NumberFormat = "@"

$ws.Range("A2").Value = "130, 455, 780"
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = "130"
$ws.Range("D2").Value = "5269"

$ws.Range("A3").Value = "1073, 1105"
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = "1105"
$ws.Range("D3").Value = "5131"

$ws.Range("A4").Value = "130, 1073, 1105"
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = "1105"
$ws.Range("D4").Value = "5399"

$ws.Range("A5").Value = "423, 748, 780, 1073"
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = "780, 780"
$ws.Range("D5").Value = "5677, 5887"

$ws.Range("A6").Value = "423, 1105, 1105"
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = "1105"
$ws.Range("D6").Value = "5331"

$ws.Range("A7").Value = "423, 1073, 1105"
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = "1105"
$ws.Range("D7").Value = "5433"

$ws.Range("A8").Value = "98, 130, 748, 1073"
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = "130"
$ws.Range("D8").Value = "5582"
